# Actualizacion lista de precios mar 16/09/2025 22:00:10,59
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista")

# Update the price values in column C (rows 1-3)
$ws.Range("C1").Value = 1003
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3

# Fix the Azucar / Sal row order + drop the accent on "Azucar"
$ws.Range("B9").Value = "Azucar x1kg"
$ws.Range("B10").Value = "Sal x1kg"

# Apply a currency number format to the whole price column (C1:C10)
$ws.Range("C1:C10").NumberFormat = '"$"\ #,##0.00'

# Column C should size itself to the new (wider) formatted contents
$ws.Columns.Item(3).ColumnWidth = 9.7

# Update the selected cell in the sheet view
$ws.Range("D4").Select()

$wb.Save()
